$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column numeric-looking text values keep their exact text representation
# (leading/trailing zeros) by forcing the cell format to Text before assignment.

$textCells = @("D2", "D3", "D4", "D5", "D6", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D40", "D41", "D44", "D45", "D47", "D48")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '247.01'
$ws.Range("D3").Value = '26.40'
$ws.Range("D4").Value = '5.086'
$ws.Range("D5").Value = '0.05618'
$ws.Range("D6").Value = '6.517'
$ws.Range("D8").Value = '0.8460'
$ws.Range("D9").Value = '0.1341'
$ws.Range("B10").Value = 'BitrueCoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D10").Value = '0.02905'
$ws.Range("E10").Value = '9BitrueCoinBTR'
$ws.Range("B11").Value = 'BitMartToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D11").Value = '0.09403'
$ws.Range("E11").Value = '10BitMartTokenBMX'
$ws.Range("B12").Value = 'BitForexToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D12").Value = '0.001521'
$ws.Range("E12").Value = '11BitForexTokenBF'
$ws.Range("B13").Value = 'One'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D13").Value = '0.0005957'
$ws.Range("E13").Value = '12OneONEWorstin24h'
$ws.Range("B14").Value = 'TigerCash'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D14").Value = '0.006130'
$ws.Range("E14").Value = '13TigerCashTCH'
$ws.Range("B15").Value = 'UpBots'
$ws.Range("C15").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D15").Value = '0.007488'
$ws.Range("E15").Value = '14UpBotsUBXTBestin24h'
$ws.Range("D16").Value = '3.589'
$ws.Range("D17").Value = '3.019'
$ws.Range("D20").Value = '0.07005'
$ws.Range("B21").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C21").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D21").Value = '0.03197'
$ws.Range("E21").Value = '20LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B22").Value = 'ProBitToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D22").Value = '0.1320'
$ws.Range("E22").Value = '21ProBitTokenPROB'
$ws.Range("B23").Value = 'MCDex'
$ws.Range("C23").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D23").Value = '3.746'
$ws.Range("E23").Value = '22MCDexMCB'
$ws.Range("B24").Value = 'CoinExToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D24").Value = '0.04660'
$ws.Range("E24").Value = '23CoinExTokenCET'
$ws.Range("B25").Value = 'ZBToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D25").Value = '0.1350'
$ws.Range("E25").Value = '24ZBTokenZB'
$ws.Range("B26").Value = 'BitKan'
$ws.Range("C26").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D26").Value = '0.001248'
$ws.Range("E26").Value = '25BitKanKAN'
$ws.Range("B27").Value = 'HotbitToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D27").Value = '0.004607'
$ws.Range("E27").Value = '26HotbitTokenHTB'
$ws.Range("B28").Value = 'NitroEx'
$ws.Range("C28").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D28").Value = '0.00009600'
$ws.Range("E28").Value = '27NitroExNTX'
$ws.Range("D40").Value = '0.03678'
$ws.Range("D41").Value = '0.006180'
$ws.Range("D44").Value = '0.008897'
$ws.Range("E44").Value = '43LocalTradersLCT'
$ws.Range("D45").Value = '0.00005289'
$ws.Range("D47").Value = '0.1499'
$ws.Range("D48").Value = '0.002533'
